$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PIM_Add_Employee")

# Update employee name data (first name, last name, derived username/password)
$ws.Range("A2").Value = "Nishchay"
$ws.Range("C2").Value = "Angra"
$ws.Range("E2").Value = "Nishchay_Angra"
$ws.Range("F2").Value = "Nishchay@20"
$ws.Range("G2").Value = "Nishchay@20"

# The confirm-password cell no longer carries its own mailto hyperlink;
# only the create-password cell (F2) keeps one. Remove just the G2 link.
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$G$2') {
        $hl.Delete()
    }
}

# Autofit column B now that it holds a wider value
$ws.Columns.Item(2).AutoFit() | Out-Null

# Reflect the final selection left by the author on this sheet
$ws.Range("G2").Select()
